# Refresh the cryptocurrency price (column D) and 1h volume change (column E)
# figures on the sheet, matching the scraper's latest run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '36.771.42'
$ws.Range("E2").Value = '  +2.87%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.130.08'
$ws.Range("E3").Value = '  +12.58%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.11%  '

# Row 5: BNB
$ws.Range("D5").Value = '''251.49'
$ws.Range("E5").Value = '  +1.71%  '

# Row 6: XRP
$ws.Range("D6").Value = '''0.673'
$ws.Range("E6").Value = '  -2.63%  '

# Row 7: USDC
$ws.Range("E7").Value = '  +0.01%  '

# Row 8: Solana
$ws.Range("D8").Value = '''45.55'
$ws.Range("E8").Value = '  +5.27%  '

# Row 9: OKB
$ws.Range("D9").Value = '''61.38'
$ws.Range("E9").Value = '  +7.67%  '

# Row 10: Cardano
$ws.Range("E10").Value = '  +2.99%  '

# Row 11: Dogecoin
$ws.Range("D11").Value = '''0.0738'
$ws.Range("E11").Value = '  -2.57%  '

# Row 12: TRON
$ws.Range("E12").Value = '  +0.58%  '

# Row 13: Chainlink
$ws.Range("D13").Value = '''14.67'
$ws.Range("E13").Value = '  -1.95%  '

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = '2.435.49'
$ws.Range("E14").Value = '  +12.29%  '

# Row 15: Polygon
$ws.Range("D15").Value = '''0.857'
$ws.Range("E15").Value = '  +8.57%  '

# Row 16: WrappedEther
$ws.Range("D16").Value = '2.125.25'
$ws.Range("E16").Value = '  +11.76%  '

# Row 17: Polkadot
$ws.Range("E17").Value = '  +0.73%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = '36.769.21'
$ws.Range("E18").Value = '  +2.78%  '

# Row 19: Litecoin
$ws.Range("D19").Value = '''73.59'
$ws.Range("E19").Value = '  +0.54%  '

# Row 20: ShibaInu
$ws.Range("E20").Value = '  -0.94%  '

# Row 21: BitcoinCash
$ws.Range("D21").Value = '''241.56'
$ws.Range("E21").Value = '  -2.01%  '

# Row 22: Avalanche
$ws.Range("E22").Value = '  -0.27%  '

# Row 23: Uniswap
$ws.Range("D23").Value = '''5.15'
$ws.Range("E23").Value = '  -0.79%  '

# Row 24: Dai
$ws.Range("D24").Value = '''1.00'
$ws.Range("E24").Value = '  -0.02%  '

# Row 25: Toncoin
$ws.Range("E25").Value = '  -8.44%  '

# Row 26: Monero
$ws.Range("D26").Value = '''170.35'
$ws.Range("E26").Value = '  +1.73%  '

# Row 27: EthereumClassic
$ws.Range("D27").Value = '''20.93'
$ws.Range("E27").Value = '  +13.45%  '

# Row 28: Cosmos
$ws.Range("D28").Value = '''9.06'
$ws.Range("E28").Value = '  +4.25%  '

# Row 29: PancakeSwap
$ws.Range("D29").Value = '''2.02'
$ws.Range("E29").Value = '  -7.36%  '

# Row 30: Stellar
$ws.Range("D30").Value = '''0.124'
$ws.Range("E30").Value = '  -3.67%  '

# Row 31: Gas
$ws.Range("D31").Value = '''22.03'
$ws.Range("E31").Value = '  +43.83%  '

# Row 32: Filecoin
$ws.Range("D32").Value = '''4.49'
$ws.Range("E32").Value = '  +0.85%  '

# Row 33: Hedera
$ws.Range("D33").Value = '''0.0600'
$ws.Range("E33").Value = '  -1.06%  '

# Row 34: Kaspa
$ws.Range("D34").Value = '''0.0917'
$ws.Range("E34").Value = '  +17.27%  '

# Row 35: LidoDAOToken
$ws.Range("E35").Value = '  +19.68%  '

# Row 37: WEMIXToken
$ws.Range("E37").Value = '  -1.18%  '

# Row 38: InternetComputer(DFINITY)
$ws.Range("D38").Value = '''4.10'
$ws.Range("E38").Value = '  -4.31%  '

# Row 39: ImmutableX
$ws.Range("D39").Value = '''0.904'
$ws.Range("E39").Value = '  +5.14%  '

# Row 40: TrustWalletToken
$ws.Range("E40").Value = '  -8.40%  '

# Row 41: ARBITRUM
$ws.Range("E41").Value = '  +10.69%  '

# Row 42: Aave
$ws.Range("D42").Value = '''101.04'
$ws.Range("E42").Value = '  +0.89%  '

# Row 43: VeChain
$ws.Range("D43").Value = '''0.0221'
$ws.Range("E43").Value = '  -2.90%  '

# Row 44: HuobiToken
$ws.Range("D44").Value = '''2.82'
$ws.Range("E44").Value = '  +16.82%  '

# Row 45: InjectiveProtocol
$ws.Range("D45").Value = '''16.40'
$ws.Range("E45").Value = '  -3.48%  '

# Row 46: Maker
$ws.Range("D46").Value = '1.371.84'
$ws.Range("E46").Value = '  +4.03%  '

# Row 47: Cronos
$ws.Range("D47").Value = '''0.0843'
$ws.Range("E47").Value = '  +4.17%  '

# Row 48: RocketPoolETH
$ws.Range("D48").Value = '2.311.34'
$ws.Range("E48").Value = '  +11.60%  '

# Row 49: MXToken
$ws.Range("D49").Value = '''2.85'
$ws.Range("E49").Value = '  +3.42%  '

# Row 50: RenderToken
$ws.Range("E50").Value = '  -2.22%  '

# Row 51: THORChain
$ws.Range("E51").Value = '  +17.60%  '
